$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table (columns B:Link, C:Link, D:Price, E:Volume(1h))
# Values that look like plain decimal numbers are entered with a leading "'" via
# .Formula so Excel keeps them as literal text (matching the source data) instead
# of re-parsing/rounding them as floating point numbers.
$ws.Cells.Item(2, 4).Value = '36.605.55'
$ws.Cells.Item(2, 5).Value = '  -1.74%  '
$ws.Cells.Item(3, 4).Value = '2.024.91'
$ws.Cells.Item(3, 5).Value = '  +0.95%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Formula = "'235.18"
$ws.Cells.Item(5, 5).Value = '  -9.30%  '
$ws.Cells.Item(6, 4).Formula = "'0.602"
$ws.Cells.Item(6, 5).Value = '  -2.56%  '
$ws.Cells.Item(7, 5).Value = '  +0.04%  '
$ws.Cells.Item(8, 4).Formula = "'54.87"
$ws.Cells.Item(8, 5).Value = '  -2.89%  '
$ws.Cells.Item(9, 5).Value = '  -2.97%  '
$ws.Cells.Item(10, 4).Formula = "'57.77"
$ws.Cells.Item(10, 5).Value = '  +2.94%  '
$ws.Cells.Item(11, 4).Formula = "'0.0750"
$ws.Cells.Item(11, 5).Value = '  -2.82%  '
$ws.Cells.Item(12, 4).Formula = "'0.101"
$ws.Cells.Item(12, 5).Value = '  -0.53%  '
$ws.Cells.Item(13, 4).Value = '2.320.00'
$ws.Cells.Item(13, 5).Value = '  +0.72%  '
$ws.Cells.Item(14, 4).Formula = "'14.18"
$ws.Cells.Item(14, 5).Value = '  -0.26%  '
$ws.Cells.Item(15, 4).Formula = "'20.13"
$ws.Cells.Item(15, 5).Value = '  -7.44%  '
$ws.Cells.Item(16, 4).Formula = "'0.765"
$ws.Cells.Item(16, 5).Value = '  -3.51%  '
$ws.Cells.Item(17, 5).Value = '  -1.93%  '
$ws.Cells.Item(18, 4).Value = '2.021.00'
$ws.Cells.Item(18, 5).Value = '  -0.41%  '
$ws.Cells.Item(19, 4).Value = '36.497.83'
$ws.Cells.Item(19, 5).Value = '  -1.98%  '
$ws.Cells.Item(20, 4).Formula = "'67.79"
$ws.Cells.Item(20, 5).Value = '  -3.24%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0798'
$ws.Cells.Item(21, 5).Value = '  -4.34%  '
$ws.Cells.Item(22, 4).Formula = "'5.37"
$ws.Cells.Item(22, 5).Value = '  +5.37%  '
$ws.Cells.Item(23, 4).Formula = "'220.36"
$ws.Cells.Item(23, 5).Value = '  -4.99%  '
$ws.Cells.Item(24, 5).Value = '  +0.13%  '
$ws.Cells.Item(25, 5).Value = '  +1.54%  '
$ws.Cells.Item(26, 4).Formula = "'2.41"
$ws.Cells.Item(26, 5).Value = '  -6.67%  '
$ws.Cells.Item(27, 4).Formula = "'163.42"
$ws.Cells.Item(27, 5).Value = '  -0.71%  '
$ws.Cells.Item(28, 4).Formula = "'8.61"
$ws.Cells.Item(28, 5).Value = '  -3.77%  '
$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(29, 4).Formula = "'0.129"
$ws.Cells.Item(29, 5).Value = '  +0.27%  '
$ws.Cells.Item(30, 2).Value = 'ImmutableX'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(30, 4).Formula = "'1.39"
$ws.Cells.Item(30, 5).Value = '  +5.40%  '
$ws.Cells.Item(31, 4).Formula = "'18.93"
$ws.Cells.Item(31, 5).Value = '  -3.09%  '
$ws.Cells.Item(32, 5).Value = '  -2.01%  '
$ws.Cells.Item(33, 5).Value = '  -4.88%  '
$ws.Cells.Item(34, 5).Value = '  -5.47%  '
$ws.Cells.Item(35, 4).Formula = "'2.46"
$ws.Cells.Item(35, 5).Value = '  +4.42%  '
$ws.Cells.Item(36, 4).Formula = "'4.24"
$ws.Cells.Item(36, 5).Value = '  -4.82%  '
$ws.Cells.Item(38, 5).Value = '  -1.99%  '
$ws.Cells.Item(39, 4).Formula = "'3.29"
$ws.Cells.Item(39, 5).Value = '  -2.65%  '
$ws.Cells.Item(40, 5).Value = '  +4.78%  '
$ws.Cells.Item(41, 4).Formula = "'2.98"
$ws.Cells.Item(41, 5).Value = '  -2.44%  '
$ws.Cells.Item(42, 2).Value = 'Maker'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(42, 4).Value = '1.456.23'
$ws.Cells.Item(42, 5).Value = '  +1.81%  '
$ws.Cells.Item(43, 2).Value = 'Cronos'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(43, 4).Formula = "'0.0928"
$ws.Cells.Item(43, 5).Value = '  +0.47%  '
$ws.Cells.Item(44, 2).Value = 'FTXToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(44, 4).Formula = "'4.21"
$ws.Cells.Item(44, 5).Value = '  +42.31%  '
$ws.Cells.Item(45, 5).Value = '  -3.23%  '
$ws.Cells.Item(46, 5).Value = '  -5.63%  '
$ws.Cells.Item(47, 4).Formula = "'90.13"
$ws.Cells.Item(47, 5).Value = '  +0.88%  '
$ws.Cells.Item(48, 4).Formula = "'15.32"
$ws.Cells.Item(48, 5).Value = '  -1.89%  '
$ws.Cells.Item(49, 5).Value = '  -1.10%  '
$ws.Cells.Item(50, 4).Formula = "'2.88"
$ws.Cells.Item(50, 5).Value = '  -1.62%  '
$ws.Cells.Item(51, 4).Formula = "'6.88"
$ws.Cells.Item(51, 5).Value = '  -1.10%  '
